# Add benchmark data for the new "Intel Core i7 6820HK" CPU.
#
# The new sheet is built by duplicating the structurally-identical
# "i7 4770k" sheet (same layout, formulas, number formats and merged
# header cell) and then overwriting just the handful of input cells that
# differ: the CPU name, its clock speed, the four "ops" counts that vary
# per category, and the eight "Actual" benchmark results. The dependent
# Max/Efficiency formulas (columns F and H) recompute automatically.

$wb = $excel.ActiveWorkbook

# The selection on the "i7 4770k" sheet moves to E24 in the final file, so
# update it (and leave that sheet itself otherwise untouched) before we
# branch off of it.
$srcSheet = $wb.Worksheets.Item("i7 4770k")
$srcSheet.Activate() | Out-Null
$srcSheet.Range("E24").Select() | Out-Null

# Duplicate "i7 4770k" to the end of the workbook to get an exact copy of
# its formulas/styles/merged cells, then rename it for the new CPU.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$srcSheet.Copy($null, $lastSheet)
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "i7 6820HK"

# CPU label (merged B2:C2) and clock speed used by the Max formulas.
$newSheet.Range("B2").Value = "Intel Core i7 6820HK"
$newSheet.Range("E2").Value = 3.20467

# "ops" counts for the first row of each SP-128/DP-128/SP-256/DP-256 block.
$newSheet.Range("E4").Value = 32
$newSheet.Range("E9").Value = 16
$newSheet.Range("E14").Value = 64
$newSheet.Range("E19").Value = 32

# Measured "Actual" results for every row (Efficiency recomputes from these).
$newSheet.Range("G4").Value = 101.391
$newSheet.Range("G5").Value = 100.002
$newSheet.Range("G6").Value = 100.322
$newSheet.Range("G7").Value = 204.45

$newSheet.Range("G9").Value = 50.8911
$newSheet.Range("G10").Value = 50.1014
$newSheet.Range("G11").Value = 50.1672
$newSheet.Range("G12").Value = 102.113

$newSheet.Range("G14").Value = 200.589
$newSheet.Range("G15").Value = 199.187
$newSheet.Range("G16").Value = 199.86
$newSheet.Range("G17").Value = 408.248

$newSheet.Range("G19").Value = 100.424
$newSheet.Range("G20").Value = 98.087
$newSheet.Range("G21").Value = 99.9052
$newSheet.Range("G22").Value = 204.166

# Make the new sheet the active tab with its own selection at E24, matching
# the final workbook state.
$newSheet.Activate() | Out-Null
$newSheet.Range("E24").Select() | Out-Null
